# Update the "Perc Decline per Doubling" label on the PDiCECpDoC sheet to
# clarify the unit is dimensionless.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDiCECpDoC")
$ws.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"

# Leave the data sheet's selection on the value cell (B2), then hand the
# active tab back to "About" so the first sheet remains the one shown
# when the workbook is opened.
$ws.Range("B2").Select()
$wb.Worksheets.Item("About").Activate()
